$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D. This shifts the existing "Tipo" column (and its
# data) from D to E, while carrying over the header formatting (bold,
# centered, bordered) so the new column lines up with the rest of the header
# row automatically.
$ws.Range("D1").EntireColumn.Insert()

# New "MAE" metric header in the freshly inserted column.
$ws.Range("D1").Value = "MAE"

# New MAE value for the single data row.
$ws.Range("D2").Value = 0.2060008420292294

# Updated MSE / R2 metric values for the single data row.
$ws.Range("B2").Value = 0.06737435472845997
$ws.Range("C2").Value = 0.9987692300954371
